$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook's shared-string table is content-addressed and append-only:
# new literal text gets pushed onto the end of the table in the order it is
# assigned, while text that already exists keeps its original slot. To land
# the new "{... .role}" strings in the same interleaved order as the
# reference edit (i18n.role right after i18n.email, contacts[i].role right
# after contacts[i].email, etc.) we first clear out the existing text so
# every cell's value becomes "new" again, then re-enter everything -
# including the role column - strictly in row-major, left-to-right order.
$ws.Range("A1:F3").ClearContents()

# Give column G (the new "role" column) the same column formatting as
# column F before widening the data block, then have it pick up each row's
# header/body formatting the same way columns A-F already do.
$ws.Range("F1:F3").Copy()
$ws.Range("G1:G3").PasteSpecial(-4122) # xlPasteFormats

# Row 1 - i18n headers
$ws.Range("A1").Value = "{d.i18n.name}"
$ws.Range("B1").Value = "{d.i18n.address}"
$ws.Range("C1").Value = "{d.i18n.unitName}"
$ws.Range("D1").Value = "{d.i18n.unitType}"
$ws.Range("E1").Value = "{d.i18n.phone}"
$ws.Range("F1").Value = "{d.i18n.email}"
$ws.Range("G1").Value = "{d.i18n.role}"

# Row 2 - contacts[i] placeholders
$ws.Range("A2").Value = "{d.contacts[i].name}"
$ws.Range("B2").Value = "{d.contacts[i].address}"
$ws.Range("C2").Value = "{d.contacts[i].unitName}"
$ws.Range("D2").Value = "{d.contacts[I].unitType}"
$ws.Range("E2").Value = "{d.contacts[i].phone}"
$ws.Range("F2").Value = "{d.contacts[i].email}"
$ws.Range("G2").Value = "{d.contacts[i].role}"

# Row 3 - contacts[i+1] placeholders
$ws.Range("A3").Value = "{d.contacts[i+1].name}"
$ws.Range("B3").Value = "{d.contacts[i+1].address}"
$ws.Range("C3").Value = "{d.contacts[i+1].unitName}"
$ws.Range("D3").Value = "{d.contacts[I+1].unitType}"
$ws.Range("E3").Value = "{d.contacts[i+1].phone}"
$ws.Range("F3").Value = "{d.contacts[i+1].email}"
$ws.Range("G3").Value = "{d.contacts[i+1].role}"

# The old template reserved rows 4-10 as blank, pre-styled filler rows.
# They are dropped entirely now that the export only ever renders 3 rows.
$ws.Range("A4:G10").EntireRow.Delete()
